$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '42.677.02'
$ws.Range("D2").Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -1.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.300.73'
$ws.Range("D3").Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '303.86'
$ws.Range("D5").Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -2.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '98.61'
$ws.Range("D6").Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -5.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.505'
$ws.Range("D7").Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -5.37%  '
$ws.Cells.Item(8, 5).Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.501'
$ws.Range("D9").Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -5.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '34.30'
$ws.Range("D10").Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -6.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '51.60'
$ws.Range("D11").Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -2.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0787'
$ws.Range("D12").Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -3.25%  '
$ws.Cells.Item(13, 5).Value = '  +0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.71'
$ws.Range("D14").Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -4.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.663.04'
$ws.Range("D15").Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -0.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '15.61'
$ws.Range("D16").Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +3.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.308.04'
$ws.Range("D17").Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.807'
$ws.Range("D18").Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '42.656.71'
$ws.Range("D19").Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0₃0898'
$ws.Range("D20").Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '11.50'
$ws.Range("D21").Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -5.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.06'
$ws.Range("D22").Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -2.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '68.98'
$ws.Range("D23").Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +1.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '234.18'
$ws.Range("D24").Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -3.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.97'
$ws.Range("D25").Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -3.18%  '
$ws.Cells.Item(26, 5).Value = '  -3.56%  '
$ws.Cells.Item(27, 5).Value = '  -0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '25.06'
$ws.Range("D28").Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.17'
$ws.Range("D29").Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -6.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '34.51'
$ws.Range("D30").Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -7.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '9.17'
$ws.Range("D31").Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -4.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '162.59'
$ws.Range("D32").Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -3.12%  '
$ws.Cells.Item(33, 5).Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.01'
$ws.Range("D34").Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -5.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.61'
$ws.Range("D35").Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +2.91%  '
$ws.Cells.Item(36, 5).Value = '  -3.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.0712'
$ws.Range("D37").Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -4.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '16.90'
$ws.Range("D38").Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -8.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.88'
$ws.Range("D39").Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -6.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.79'
$ws.Range("D40").Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -4.70%  '
$ws.Cells.Item(41, 5).Value = '  -5.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.111'
$ws.Range("D42").Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -4.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.47'
$ws.Range("D43").Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -9.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.986.54'
$ws.Range("D44").Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '18.69'
$ws.Range("D45").Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -2.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0279'
$ws.Range("D46").Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -4.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '10.22'
$ws.Range("D47").Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +2.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.86'
$ws.Range("D48").Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -6.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '55.24'
$ws.Range("D49").Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -1.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.85'
$ws.Range("D50").Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -3.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.528.97'
$ws.Range("D51").Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.35%  '
